# Adding the changes we made on may 9th
# Appends 10 new rows (timestamps 2000-2900) of walkingToRunning sensor
# data below the existing data on Sheet1 (rows 22-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2000, "walkingToRunning", -40.18527779212365,-12.29963522690993,-2.489131450653076,-4.548879146575928,1.983122944831848,-2.030394077301025),
    @(2100, "walkingToRunning", -5.963009345225885,-38.13061856000844,21.80784027393036,0.2330931723117828,7.29985761642456,-8.993494987487793),
    @(2200, "walkingToRunning", 6.794984866411211,-17.97448307428637,8.985512146583247,4.196409702301025,-7.729958534240723,-5.114867210388184),
    @(2300, "walkingToRunning", 8.772580513587364,-2.704863694997934,-3.532207745772141,1.506084203720093,-3.158908605575561,-2.592853784561157),
    @(2400, "walkingToRunning", 15.409985371125,-21.02941689124464,-3.300734495505285,0.7066032290458679,-5.12858247756958,0.3414172530174255),
    @(2500, "walkingToRunning", 29.99191743899604,-16.15106448149102,-5.187779292082201,-3.361641645431519,6.300107002258301,2.977946758270264),
    @(2600, "walkingToRunning", 43.31018800001878,-13.67232219989484,-2.019487449756035,-2.117412805557251,13.40914916992188,-1.111603498458862),
    @(2700, "walkingToRunning", 47.68574269612633,-58.8819274902343,7.092501203219018,1.100218415260315,1.102149248123169,-2.866626739501953),
    @(2800, "walkingToRunning", 27.69281722337806,-49.70874527173221,14.84765016115646,0.4749082028865814,-12.12523937225342,-7.235276222229004),
    @(2900, "walkingToRunning", -3.377501120934119,4.849835469172549,0.1248559951782222,3.548196077346802,9.086838722229004,1.108407616615295),
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
